$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G (header "K") values for rows 2-5
$ws.Range("G2").Value = 4
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
